# Add season record columns (Wins, Losses, Ties) to the PIT_1997 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy formatting from the existing header style (A1) so the new
# header cells (AD1:AF1) match the bold/centered/bordered look of the rest
# of row 1, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record data: every player row (2-45) gets the team's season
# record of 79 wins, 83 losses, 0 ties.
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
